# Generate Report for Handback
#
# Refresh the generated handback-status report: the second source file
# (f87a9faf-9bbb-4e63-b66a-25961983642a.md) has gone through another
# handoff/handback cycle, so its timestamps (and the roll-up "Latest HO
# Xliff Generate Date" on the Overview sheet) need to be brought forward.
#
# Row 2 (9e5a7da1-...) is untouched; only row 3 (f87a9faf-...) changes,
# on the "zh-cn" sheet, the "de-de" sheet, and the roll-up on "Overview".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# zh-cn: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K)
$zhcn.Range("H3").Value = "2016-09-05 04:54:08"
$zhcn.Range("K3").Value = "2016-09-05 04:54:26"

# de-de: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K)
$dede.Range("H3").Value = "2016-09-05 04:54:14"
$dede.Range("K3").Value = "2016-09-05 04:54:33"

# Overview: Latest HO Xliff Generate Date (G) for the same file
$overview.Range("G3").Value = "2016-09-05 04:54:14"
